$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark existing trade people as Default
$ws.Range("B9").Value = "Roger Yearwood(Default)"
$ws.Range("B10").Value = "Ralph Carrow Default)"
$ws.Range("C9").Value = "Elwood Olin Default)"

# New contractor entities
$ws.Range("B11").Value = "Sonny Getchell (Contractor)"
$ws.Range("B12").Value = "Taylor Diniz((Contractor))"
$ws.Range("C10").Value = "Taylor Diniz((Contractor))"

# New property row
$ws.Range("B5").Value = "Simon's Property"

# New allocation notes
$ws.Range("B14").Value = "Sonny's Allocation: Simon's "

$ws.Range("I5").Value = "Share same property as Lee"
$ws.Range("I6").Value = "yes"

$ws.Range("C14").Value = "Taylor's Allocation: Lee's"
$ws.Range("B15").Value = "Taylor's Allocation: Joe's"

# Update selection to match final cursor position
$ws.Range("C15").Select()
